$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings (e.g. "1.000", "241.72") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Data rows: row index, Coin, Link, Price, Volume(1h)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.040.70', '  -0.68%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.833.24', '  -0.61%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9988', '  -0.02%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '241.72', '  -0.40%  '),
    @(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6285', '  -5.36%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  +0.05%  '),
    @(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07612', '  +2.05%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2918', '  -1.39%  '),
    @(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '22.73', '  -3.11%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07743', '  -0.47%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.828.42', '  -0.86%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.961', '  -1.30%  '),
    @(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6651', '  -1.33%  '),
    @(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '82.82', '  -1.08%  '),
    @(16, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009776', '  +13.61%  '),
    @(17, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.006', '  -3.03%  '),
    @(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.024.98', '  -0.85%  '),
    @(19, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '225.53', '  -1.14%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.34', '  -1.85%  '),
    @(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9999', '  -0.04%  '),
    @(22, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.198', '  -0.38%  '),
    @(23, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.04%  '),
    @(24, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '158.20', '  -0.53%  '),
    @(25, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.428', '  -2.60%  '),
    @(26, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1372', '  -2.80%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.88', '  -1.19%  '),
    @(28, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.491', '  -1.44%  '),
    @(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.071', '  -1.69%  '),
    @(30, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.022', '  -0.94%  '),
    @(31, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.198', '  +0.51%  '),
    @(32, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05193', '  -2.77%  '),
    @(33, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.847', '  -1.83%  '),
    @(34, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7397', '  -1.42%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.144', '  -1.24%  '),
    @(36, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.699', '  +1.75%  '),
    @(37, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.267.23', '  -3.85%  '),
    @(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.760', '  +0.16%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01788', '  -0.82%  '),
    @(40, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.275', '  -2.26%  '),
    @(41, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8950', '  -2.26%  '),
    @(42, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  +0.11%  '),
    @(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.51', '  -1.86%  '),
    @(44, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '1.979.19', '  -2.65%  '),
    @(45, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000124', '  +1.74%  '),
    @(46, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '64.68', '  -2.16%  '),
    @(47, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.5113', '  -0.58%  '),
    @(48, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3997', '  -0.91%  '),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '8.845', '  +0.78%  '),
    @(50, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05755', '  -1.88%  '),
    @(51, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.632', '  -7.21%  '),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}
